$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 8-15 get new B-column labels ("line7"/"line8" pushed in ahead of the
# "extr*" labels, which shift down by two) plus new C/D/E data values.
# Rows 16-17 are brand new rows appended at the bottom.

$rows = @(
    @{ Row = 8;  A = 6;  B = "line7"; C = 14; D = 11; E = $false },
    @{ Row = 9;  A = 7;  B = "line8"; C = 16; D = 9;  E = $true  },
    @{ Row = 10; A = 8;  B = "extr1"; C = 5;  D = 12; E = $true  },
    @{ Row = 11; A = 9;  B = "extr2"; C = 5;  D = 9;  E = $true  },
    @{ Row = 12; A = 10; B = "extr3"; C = 10; D = 11; E = $true  },
    @{ Row = 13; A = 11; B = "extr4"; C = 7;  D = 8;  E = $true  },
    @{ Row = 14; A = 12; B = "extr5"; C = 9;  D = 11; E = $false },
    @{ Row = 15; A = 13; B = "extr6"; C = 7;  D = 11; E = $true  },
    @{ Row = 16; A = 14; B = "extr7"; C = 5;  D = 7;  E = $false },
    @{ Row = 17; A = 15; B = "extr8"; C = 8;  D = 5;  E = $true  }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}

# New rows 16 & 17 need the same formatting (bold/centered/bordered) that
# column A already carries on every other data row (copy it from A15).
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false
